# Apply updated crypto price/volume figures to the sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $val) {
    # Force the cell to retain an exact text representation (e.g. trailing
    # zeros like "8.60") instead of Excel auto-coercing it to a number.
    $r = $ws.Range($cell)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.ClearFormats()
}

$ws.Range('D2').Value = '60.032.82'
$ws.Range('D3').Value = '2.418.01'
$ws.Range('E3').Value = '  -1.17%  '
$ws.Range('E4').Value = '  -0.01%  '
Set-TextValue 'D5' '552.30'
$ws.Range('E5').Value = '  -0.81%  '
Set-TextValue 'D6' '137.26'
$ws.Range('E6').Value = '  -1.20%  '
$ws.Range('E7').Value = '  +0.01%  '
$ws.Range('E8').Value = '  +4.09%  '
$ws.Range('E9').Value = '  -1.91%  '
$ws.Range('E10').Value = '  -2.47%  '
$ws.Range('E11').Value = '  -1.10%  '
Set-TextValue 'D12' '0.353'
$ws.Range('E12').Value = '  -2.16%  '
$ws.Range('E13').Value = '  +1.68%  '
$ws.Range('D14').Value = '2.847.39'
$ws.Range('E14').Value = '  -1.25%  '
$ws.Range('D15').Value = '59.957.34'
$ws.Range('E15').Value = '  -0.26%  '
$ws.Range('E16').Value = '  -2.46%  '
$ws.Range('D17').Value = '2.412.10'
$ws.Range('E17').Value = '  -1.78%  '
$ws.Range('E18').Value = '  -1.73%  '
$ws.Range('E19').Value = '  -0.68%  '
Set-TextValue 'D20' '329.17'
$ws.Range('E20').Value = '  -1.85%  '
$ws.Range('E21').Value = '  -4.03%  '
$ws.Range('E22').Value = '  +0.05%  '
Set-TextValue 'D23' '65.68'
$ws.Range('E23').Value = '  +1.25%  '
$ws.Range('E24').Value = '  +4.22%  '
Set-TextValue 'D25' '8.60'
$ws.Range('E25').Value = '  +0.18%  '
$ws.Range('E26').Value = '  +0.11%  '
$ws.Range('E27').Value = '  -0.59%  '
$ws.Range('E28').Value = '  -2.80%  '
$ws.Range('E29').Value = '  -2.51%  '
Set-TextValue 'D30' '168.84'
$ws.Range('E30').Value = '  -1.33%  '
Set-TextValue 'D31' '6.06'
$ws.Range('E31').Value = '  -4.29%  '
Set-TextValue 'D32' '18.58'
$ws.Range('E32').Value = '  -1.41%  '
$ws.Range('E33').Value = '  -0.48%  '
$ws.Range('E35').Value = '  -0.65%  '
$ws.Range('E36').Value = '  +0.06%  '
$ws.Range('E37').Value = '  -2.73%  '
$ws.Range('E38').Value = '  -2.63%  '
Set-TextValue 'D39' '321.07'
$ws.Range('E39').Value = '  +1.29%  '
$ws.Range('E40').Value = '  -3.05%  '
$ws.Range('E41').Value = '  -2.15%  '
Set-TextValue 'D42' '139.97'
$ws.Range('E42').Value = '  -2.88%  '
$ws.Range('E43').Value = '  +0.68%  '
Set-TextValue 'D44' '19.59'
$ws.Range('E44').Value = '  +0.62%  '
Set-TextValue 'D45' '0.0515'
$ws.Range('E45').Value = '  -2.23%  '
Set-TextValue 'D46' '0.577'
$ws.Range('E46').Value = '  +0.30%  '
$ws.Range('E47').Value = '  -1.75%  '
Set-TextValue 'D48' '0.387'
$ws.Range('E48').Value = '  -5.48%  '
$ws.Range('E49').Value = '  +0.02%  '
$ws.Range('E50').Value = '  -3.77%  '
$ws.Range('E51').Value = '  -1.06%  '
